# Splits the last bullet item of the list ("Algunas entidades que tienen
# datos autonumericos ... reversible también)") into three paragraphs:
#   1) the original text (unchanged);
#   2) a new bullet about pagination/sorting ("Con respecto a la
#      paginación y ordenación ...");
#   3) an (empty) bullet that now carries the trailing _GoBack bookmark.

$d = $word.ActiveDocument

# Step 1: split the paragraph in two places using Find/Replace with the
# special "^p" paragraph-mark code. Executing the split through Find
# naturally re-homes the (zero-width) _GoBack bookmark onto the new,
# now-last paragraph, exactly as real Word does when you press Enter
# right before a bookmark sitting at the end of a paragraph.
$find = $d.Content.Find
$anchor = "de manera reversible también)"
$placeholder = "PLACEHOLDER_NEW_PARAGRAPH"
$replacement = "de manera reversible también)^p" + $placeholder + "^p"
$found = $find.Execute($anchor, $true, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)
if (-not $found) {
    throw "Could not locate the anchor paragraph text to split."
}

# Step 2: locate the freshly created placeholder paragraph (the middle of
# the three) and replace its contents with the real run/proofErr
# structure that mirrors how Word marks words it doesn't recognise while
# spell-checking.
$paraStart = -1
$paraEnd = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*$placeholder*") {
        $paraStart = $p.Range.Start
        $paraEnd = $p.Range.End - 1
        break
    }
}
if ($paraStart -lt 0) {
    throw "Could not locate the placeholder paragraph."
}
$contentRange = $d.Range($paraStart, $paraEnd)

$body = @'
<w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve">Con respecto a la paginación y ordenación tenemos lo siguiente: debería hacerse en el servidor y no en el proyecto de .net (para efectos del ejercicio se </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>realizo</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> en el proyecto .net), es discutible si esto es de la capa de negocios BLL o es mas de otra capa la coloque en el mediator para no colocarle esta responsabilidad al BLL, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> pendiente encontrar una solución genérica para el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>sortby</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> y el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>sort</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t>direction</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-CO"/></w:rPr><w:t xml:space="preserve"> (no los implemente por tiempo)</w:t></w:r>
'@

$xmlSnippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + '<w:p>' + $body + '</w:p>' + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$contentRange.InsertXML($xmlSnippet)
